# Update CTI slope comparison table with 202051023 classification model run.
# Replace the estimate / std.error / df / statistic / p.value values for each
# of the three contrast rows (Deep - Mid, Deep - Shallow, Mid - Shallow).
#
# Each cell is addressed explicitly via Table.Cell(row, col) and the text is
# swapped with a Find/Replace scoped to that cell's character range, so the
# existing run formatting (Calibri, sz 20) is preserved and no stray matches
# in other cells (e.g. the repeated "108" / "0.028" values) are touched.
#
# Note: Cell.Range itself is not reliably used as the Find anchor in this
# runtime (Find.Execute on it can match content outside the cell), so the
# start/end offsets are read off Cell.Range and a fresh Document.Range(start,
# end) is built from them to scope the search precisely.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellValue($table, $row, $col, $oldText, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $d.Range($cell.Range.Start, $cell.Range.End)
    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, `
                       $true, 0, $false, $newText, 2) | Out-Null
}

# Row 2: Deep - Mid
Set-CellValue $t 2 2 "0.000" "0.007"
Set-CellValue $t 2 3 "0.019" "0.018"
Set-CellValue $t 2 4 "108"   "153"
Set-CellValue $t 2 5 "0.020" "0.375"
Set-CellValue $t 2 6 "0.984" "0.708"

# Row 3: Deep - Shallow
Set-CellValue $t 3 2 "0.028" "-0.007"
Set-CellValue $t 3 3 "0.021" "0.019"
Set-CellValue $t 3 4 "108"   "153"
Set-CellValue $t 3 5 "1.350" "-0.359"
Set-CellValue $t 3 6 "0.180" "0.720"

# Row 4: Mid - Shallow
Set-CellValue $t 4 2 "0.028" "-0.014"
Set-CellValue $t 4 3 "0.020" "0.015"
Set-CellValue $t 4 4 "108"   "153"
Set-CellValue $t 4 5 "1.376" "-0.886"
Set-CellValue $t 4 6 "0.172" "0.377"

Write-Output "Updated table with 202051023 classification model run"
